# Fix SQL error caused by an apostrophe/quote character appearing in an
# imported value ("Streep" -> "Streep'"), per the RosarioSIS commit fixing
# CSV-import SQL errors when a name contains a quote. Version 1.8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Streep" + [char]0x2019

# Leave the active cell parked on the edited cell (matches the saved
# selection state after making the edit).
$ws.Range("B3").Select()
